# add tag for block
# Insert a new "Tag" row right after the header row on the "Property" sheet.
# Existing data rows (2-13) shift down by one (to 3-14); a new row 2 is
# populated with the Tag field definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 13

# Shift existing data rows down by one, bottom-up so we never clobber a row
# before it has been copied.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":L" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":L" + $dstRow)
    $src.Copy($dst)
}

# Populate the freed-up row 2 with the new "Tag" field, following the same
# pattern as every other field row (string type, Public/Private/Save false,
# View true, Index/SaveInterval 0, RelationValue "Friend").
$ws.Range("A2").Value = "Tag"
$ws.Range("B2").Value = "string"
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = $false
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = $true
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "Friend"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""

# Match the author's resulting selection.
$ws.Range("A3").Select()
